# Rename the "aquisicoes" sheet to "Folha1", make it the active/selected
# sheet (tab), and update its selected cell to G21.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("aquisicoes")
$ws.Name = "Folha1"

# Make this sheet the active tab (was Sheet1 before).
$ws.Activate()

# Update the selection on the renamed sheet.
$ws.Range("G21").Select()
